$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PR Document")

# Fill in row 9 with the next purchase order record (mirrors the pattern of row 8)
# Copy the number formats from row 8 first so the new values keep the
# same look (text-formatted PO number, date-formatted order date, etc.)
$ws.Range("A8:G8").Copy() | Out-Null
$ws.Range("A9:G9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A9").Value = "Dept002"
$ws.Range("B9").Value = "PO_10001008"
$ws.Range("C9").Value = 44579
$ws.Range("D9").Value = "Equipment AB"
$ws.Range("E9").Value = "Equipment "
$ws.Range("F9").Value = "V0002"
$ws.Range("G9").Value = 6000

# Update the active selection to match the saved view state
$ws.Range("D10").Select()
